# ------------------------------------------------------------------
# Rewrites the "Input" sheet (sheet1) to the new standard template
# layout (16 columns A:P instead of 17 columns A:Q), and drops the
# now-empty trailing "비고" placeholder cells on "갑지"/"을지"
# (sheet2 / sheet3).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------- Sheet 1 : Input -----------------------------------
$ws1 = $wb.Worksheets.Item("Input")

# The new template only has 16 columns (A:P) - drop the old trailing
# "Q" column (old "비고") completely, cells and all.
$ws1.Range("Q1").EntireColumn.Delete()

# Wipe every value currently in A1:P4 (old 17-col layout minus Q).
$ws1.Range("A1:P4").ClearContents()

# The new header row is NOT bold/centered/bordered any more - reset
# it back to the default ("Normal") style.
$ws1.Range("A1:P1").Style = "Normal"

# New header row (A1:P1).
$headers = @("발주일자","납기일자","거래처명","거래처 이메일","납품처명","납품처 이메일","프로젝트명","대분류","중분류","소분류","품목명","규격","수량","단가","총금액","비고")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws1.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Date-looking text (발주일자/납기일자) must stay literal text, not get
# auto-parsed into a date serial by Excel - force a text number format
# before assigning, then drop back to the default style so no stray
# format survives on the saved cell.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws1.Cells.Item(2, 1) "2025-09-17"
Set-TextValue $ws1.Cells.Item(2, 2) "2025-10-08"
$ws1.Cells.Item(2, 3).Value = "유니모터스"
$ws1.Cells.Item(2, 4).Value = "유니모터스@example.com"
$ws1.Cells.Item(2, 5).Value = "힐스테이트 도곡동1차"
$ws1.Cells.Item(2, 6).Value = "delivery@example.com"
$ws1.Cells.Item(2, 7).Value = "힐스테이트 도곡동1차"
$ws1.Cells.Item(2, 8).Value = "4. 장비비"
$ws1.Cells.Item(2, 9).Value = "1) 장비비"
$ws1.Cells.Item(2, 10).Value = "기타"
$ws1.Cells.Item(2, 11).Value = "자재하차"
$ws1.Cells.Item(2, 12).Value = "KS규격-1"
$ws1.Cells.Item(2, 13).Value = 2
$ws1.Cells.Item(2, 14).Value = 80000
$ws1.Cells.Item(2, 15).Value = 176000

# Row 3
Set-TextValue $ws1.Cells.Item(3, 1) "2025-09-06"
Set-TextValue $ws1.Cells.Item(3, 2) "2025-09-18"
$ws1.Cells.Item(3, 3).Value = "유니모터스"
$ws1.Cells.Item(3, 4).Value = "유니모터스@example.com"
$ws1.Cells.Item(3, 5).Value = "힐스테이트 도곡동1차"
$ws1.Cells.Item(3, 6).Value = "delivery@example.com"
$ws1.Cells.Item(3, 7).Value = "힐스테이트 도곡동1차"
$ws1.Cells.Item(3, 8).Value = "4. 장비비"
$ws1.Cells.Item(3, 9).Value = "1) 장비비"
$ws1.Cells.Item(3, 10).Value = "기타"
$ws1.Cells.Item(3, 11).Value = "자재하차"
$ws1.Cells.Item(3, 12).Value = "KS규격-2"
$ws1.Cells.Item(3, 13).Value = 2
$ws1.Cells.Item(3, 14).Value = 80000
$ws1.Cells.Item(3, 15).Value = 176000

# Row 4
Set-TextValue $ws1.Cells.Item(4, 1) "2025-09-14"
Set-TextValue $ws1.Cells.Item(4, 2) "2025-09-05"
$ws1.Cells.Item(4, 3).Value = "유니모터스"
$ws1.Cells.Item(4, 4).Value = "유니모터스@example.com"
$ws1.Cells.Item(4, 5).Value = "힐스테이트 도곡동1차"
$ws1.Cells.Item(4, 6).Value = "delivery@example.com"
$ws1.Cells.Item(4, 7).Value = "힐스테이트 도곡동1차"
$ws1.Cells.Item(4, 8).Value = "4. 장비비"
$ws1.Cells.Item(4, 9).Value = "1) 장비비"
$ws1.Cells.Item(4, 10).Value = "기타"
$ws1.Cells.Item(4, 11).Value = "자재하차"
$ws1.Cells.Item(4, 12).Value = "KS규격-3"
$ws1.Cells.Item(4, 13).Value = 1
$ws1.Cells.Item(4, 14).Value = 80000
$ws1.Cells.Item(4, 15).Value = 88000

# ---------------- Sheet 2 ("갑지") & Sheet 3 ("을지") ----------------
# Both sheets already carry a blank placeholder "비고" inline-string
# cell in I2:I4 with no content - drop those empty cells entirely.
foreach ($name in "갑지", "을지") {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("I2").ClearContents()
    $ws.Range("I3").ClearContents()
    $ws.Range("I4").ClearContents()
}
